# The document currently splits the first sentence's leading text and the
# trailing text around the "inline   code" run into extra runs that only
# hold a single space (an artifact of earlier edits). The target layout
# merges each of those stray space-only runs into its neighbour so that:
#   "This is an example of" + " "              -> "This is an example of "
#   " " + "with three spaces."                 -> " with three spaces."
#
# Re-running Find/Replace across each run boundary with the combined text
# (old text == new text) is enough to make Word re-flow/merge the backing
# runs into a single run, matching the target OOXML exactly.

$d = $word.ActiveDocument

$d.Content.Find.Execute("This is an example of ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "This is an example of ", 2)

$d.Content.Find.Execute(" with three spaces.", $true, $false, $false, `
    $false, $false, $true, 1, $false, " with three spaces.", 2)
